# Reorder the "Estado de Cuenta" detail table (rows 16-36) so that the data is
# grouped by Periodo Mora (ascending 2501 -> 2507) instead of by worker, with
# each period block listing the three workers in the same fixed order
# (LUIS ANGEL LUNA ESCORCIA, SEBASTIAN DE JESUS RUIZ AVILA, JAIRO JESUS
# CABARCAS MARTINEZ). The per-worker Valor Mora / Salario Basico figures are
# carried along unchanged with their owning worker+period pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New row order, each entry is:
# TipoDoc, NumDoc, Nombre, Periodo, ValorMora, SalarioBasico
$data = @(
    ,@("CC", "1143401756", "LUIS ANGEL LUNA ESCORCIA",      "2501", 128000, 3200000)
    ,@("CC", "1007230501", "SEBASTIAN DE JESUS RUIZ AVILA", "2501", 76000,  1900000)
    ,@("CC", "1143351761", "JAIRO JESUS CABARCAS MARTINEZ", "2501", 80000,  2000000)
    ,@("CC", "1143401756", "LUIS ANGEL LUNA ESCORCIA",      "2502", 128000, 3200000)
    ,@("CC", "1007230501", "SEBASTIAN DE JESUS RUIZ AVILA", "2502", 76000,  1900000)
    ,@("CC", "1143351761", "JAIRO JESUS CABARCAS MARTINEZ", "2502", 80000,  2000000)
    ,@("CC", "1143401756", "LUIS ANGEL LUNA ESCORCIA",      "2503", 128000, 3200000)
    ,@("CC", "1007230501", "SEBASTIAN DE JESUS RUIZ AVILA", "2503", 76000,  1900000)
    ,@("CC", "1143351761", "JAIRO JESUS CABARCAS MARTINEZ", "2503", 80000,  2000000)
    ,@("CC", "1143401756", "LUIS ANGEL LUNA ESCORCIA",      "2504", 128000, 3200000)
    ,@("CC", "1007230501", "SEBASTIAN DE JESUS RUIZ AVILA", "2504", 76000,  1900000)
    ,@("CC", "1143351761", "JAIRO JESUS CABARCAS MARTINEZ", "2504", 80000,  2000000)
    ,@("CC", "1143401756", "LUIS ANGEL LUNA ESCORCIA",      "2505", 128000, 3200000)
    ,@("CC", "1007230501", "SEBASTIAN DE JESUS RUIZ AVILA", "2505", 76000,  1900000)
    ,@("CC", "1143351761", "JAIRO JESUS CABARCAS MARTINEZ", "2505", 80000,  2000000)
    ,@("CC", "1143401756", "LUIS ANGEL LUNA ESCORCIA",      "2506", 128000, 3200000)
    ,@("CC", "1007230501", "SEBASTIAN DE JESUS RUIZ AVILA", "2506", 76000,  1900000)
    ,@("CC", "1143351761", "JAIRO JESUS CABARCAS MARTINEZ", "2506", 80000,  2000000)
    ,@("CC", "1143401756", "LUIS ANGEL LUNA ESCORCIA",      "2507", 119467, 3200000)
    ,@("CC", "1007230501", "SEBASTIAN DE JESUS RUIZ AVILA", "2507", 52267,  1900000)
    ,@("CC", "1143351761", "JAIRO JESUS CABARCAS MARTINEZ", "2507", 74667,  2000000)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]

    $ws.Cells.Item($row, 2).Value = $rec[0]  # B: Tipo Doc Trabajador
    $ws.Cells.Item($row, 3).Value = $rec[1]  # C: N Doc Trabajador
    $ws.Cells.Item($row, 4).Value = $rec[2]  # D: Nombre Trabajador
    $ws.Cells.Item($row, 5).Value = $rec[3]  # E: Periodo Mora (text "2501".."2507")
    $ws.Cells.Item($row, 6).Value = $rec[4]  # F: Valor Mora
    $ws.Cells.Item($row, 7).Value = $rec[5]  # G: Salario Basico
}
